$wb = $excel.ActiveWorkbook

# 1. Create the new "Croatia" sheet by copying "Turkey" (closest template: same
#    columns/merges/styles), placing it after Turkey (i.e. at the end).
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Copy($null, $turkey)
$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# 2. Croatia has one extra accessory row versus Turkey ("PR1D2-Unmonitored").
#    Insert a fresh row 9 (pushes RDS800/Wg/Accessories down to 10/11/12) and
#    copy formatting down from row 8 so the new row matches the others.
$croatia.Rows.Item(9).Insert()
$croatia.Range("A8").Copy()
$croatia.Range("A9").PasteSpecial(-4122)
$croatia.Range("A9").Value = "PR1D2-Unmonitored"

# 3. Market-specific values for Croatia. (Order matters for shared-string
#    table placement: NGC-3139/T2494/T2485 lands before Croatia Market.)
$croatia.Range("B4").Value = "NGC-3139/T2494/T2485"
$croatia.Range("B2").Value = "Croatia Market"

# 4. Selection state ends on A9 for the new sheet.
$croatia.Range("A9").Select()

# 5. Turkey is no longer the active tab, so it loses its highlighted
#    selection and reverts to a whole-sheet selection like the other
#    inactive sheets.
$turkey.Range("A1:XFD1048576").Select()

# 6. Incidental selection change left on Portugal from the editing session.
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Range("A10").Select()

# 7. Leave Croatia as the active sheet/tab.
$croatia.Activate()
$croatia.Range("A9").Select()
